$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("C2").Value = "北京·EXA·全职高手ONLY·夏令营"
$ws.Range("D2").Value = "金盏路6号 蓝可可亲子乐园(金盏店)"
$ws.Range("E2").Value = "2024.06.29 10:00-06.29 16:00"
$ws.Range("F2").Value = 560
$ws.Range("G2").Value = 98
$ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=83977"
$ws.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202404/JykdQ3eR1712654574985.jpeg"
$ws.Range("F4").Value = 96
$ws.Range("F5").Value = 919
$ws.Range("F6").Value = 62
$ws.Range("F7").Value = 7172
$ws.Range("F8").Value = 113
$ws.Range("F9").Value = 159
$ws.Range("F10").Value = 6564
$ws.Range("F13").Value = 4538
$ws.Range("F16").Value = 55
$ws.Range("F17").Value = 4596
$ws.Range("F18").Value = 16
$ws.Range("F28").Value = 8234
$ws.Range("F30").Value = 1428
$ws.Range("F32").Value = 726
$ws.Range("F37").Value = 1690
$ws.Range("F41").Value = 4272
$ws.Range("F42").Value = 357
$ws.Range("F43").Value = 632
$ws.Range("F44").Value = 120
$ws.Range("F46").Value = 853
$ws.Range("F49").Value = 27

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 22

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("C5").Value = "北京·EXA·全职高手ONLY·夏令营"
$ws.Range("D5").Value = "金盏路6号 蓝可可亲子乐园(金盏店)"
$ws.Range("E5").Value = "2024.06.29 10:00-06.29 16:00"
$ws.Range("F5").Value = 560
$ws.Range("G5").Value = 98
$ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=83977"
$ws.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202404/JykdQ3eR1712654574985.jpeg"
$ws.Range("F7").Value = 96
$ws.Range("F8").Value = 22
$ws.Range("F9").Value = 919
$ws.Range("F10").Value = 62
$ws.Range("F11").Value = 7173
$ws.Range("F12").Value = 113
$ws.Range("F13").Value = 159
$ws.Range("F14").Value = 6564
$ws.Range("F17").Value = 4538
$ws.Range("F20").Value = 55
$ws.Range("F21").Value = 4596
$ws.Range("F22").Value = 16
$ws.Range("F29").Value = 8234
$ws.Range("F31").Value = 1428
$ws.Range("F33").Value = 726
$ws.Range("F37").Value = 1690
$ws.Range("F41").Value = 4273
$ws.Range("F42").Value = 357
$ws.Range("F43").Value = 632
$ws.Range("F44").Value = 120
$ws.Range("F46").Value = 853
$ws.Range("F49").Value = 27
